$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 38
$ws.Range("H38").Value = 2150917.8
$ws.Range("I38").Value = 4032378.5
$ws.Range("J38").Value = 676.5714
$ws.Range("K38").Value = 12097135.5
$ws.Range("L38").Value = 2029.7142
$ws.Range("M38").Value = -12096763.5
$ws.Range("N38").Value = -2773.7142
# row 112
$ws.Range("H112").Value = 1018.7143
$ws.Range("J112").Value = 1030.1666
$ws.Range("L112").Value = 3090.4998
$ws.Range("N112").Value = -5306.4998
# row 113
$ws.Range("H113").Value = 144886.42
$ws.Range("I113").Value = 252051.25
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 252051.25
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -248797.25
$ws.Range("N113").Value = -8508
# row 116
$ws.Range("H116").Value = 1985.9412
$ws.Range("I116").Value = 1579
$ws.Range("J116").Value = 2155.5
$ws.Range("K116").Value = 1579
$ws.Range("L116").Value = 2155.5
$ws.Range("M116").Value = 1863
$ws.Range("N116").Value = -9039.5
# row 129
$ws.Range("H129").Value = 2871.2156
$ws.Range("I129").Value = 13022.875
$ws.Range("J129").Value = 982.5349
$ws.Range("K129").Value = 39068.625
$ws.Range("L129").Value = 2947.6047
$ws.Range("M129").Value = -34068.625
$ws.Range("N129").Value = -12947.6047

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 37308.75
$ws.Range("I2").Value = 1305.8125
$ws.Range("J2").Value = 85312.664
$ws.Range("K2").Value = 1305.8125
$ws.Range("L2").Value = 85312.664
$ws.Range("M2").Value = -1192.8125
$ws.Range("N2").Value = -85538.664
# row 32
$ws.Range("H32").Value = 22816.6
$ws.Range("I32").Value = 3803.1343
$ws.Range("K32").Value = 3803.1343
$ws.Range("M32").Value = -3516.1343
# row 45
$ws.Range("H45").Value = 1775.6538
$ws.Range("I45").Value = 1681.1875
$ws.Range("J45").Value = 1926.8
$ws.Range("K45").Value = 1681.1875
$ws.Range("L45").Value = 1926.8
$ws.Range("M45").Value = -1304.1875
$ws.Range("N45").Value = -2680.8
# row 116
$ws.Range("H116").Value = 37308.75
$ws.Range("I116").Value = 1305.8125
$ws.Range("J116").Value = 85312.664
$ws.Range("K116").Value = 1305.8125
$ws.Range("L116").Value = 85312.664
$ws.Range("M116").Value = 988.1875
$ws.Range("N116").Value = -89900.664

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 37308.75
$ws.Range("I3").Value = 1305.8125
$ws.Range("J3").Value = 85312.664
$ws.Range("K3").Value = 1305.8125
$ws.Range("L3").Value = 85312.664
$ws.Range("M3").Value = -1191.8125
$ws.Range("N3").Value = -85540.664
# row 52
$ws.Range("H52").Value = 43407.875
$ws.Range("I52").Value = 44894.715
$ws.Range("J52").Value = 33000
$ws.Range("K52").Value = 44894.715
$ws.Range("L52").Value = 33000
$ws.Range("M52").Value = -44631.715
$ws.Range("N52").Value = -33526
# row 105
$ws.Range("H105").Value = 70806.93
$ws.Range("I105").Value = 45190.87
$ws.Range("K105").Value = 45190.87
$ws.Range("M105").Value = -43443.87
# row 116
$ws.Range("H116").Value = 47000
$ws.Range("J116").Value = 47000
$ws.Range("L116").Value = 47000
$ws.Range("N116").Value = -56178
# row 118
$ws.Range("H118").Value = 26800
$ws.Range("J118").Value = 26800
$ws.Range("L118").Value = 26800
$ws.Range("N118").Value = -30114
# row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# row 121
$ws.Range("H121").Value = 43407.875
$ws.Range("I121").Value = 44894.715
$ws.Range("J121").Value = 33000
$ws.Range("K121").Value = 44894.715
$ws.Range("L121").Value = 33000
$ws.Range("M121").Value = -43147.715
$ws.Range("N121").Value = -36494

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 741.5
$ws.Range("I22").Value = 386.4
$ws.Range("J22").Value = 1333.3334
$ws.Range("K22").Value = 386.4
$ws.Range("L22").Value = 1333.3334
$ws.Range("M22").Value = -36.39999999999998
$ws.Range("N22").Value = -2033.3334
# row 99
$ws.Range("H99").Value = 8894.111
$ws.Range("I99").Value = 3426.182
$ws.Range("K99").Value = 3426.182
$ws.Range("M99").Value = -1928.182
# row 109
$ws.Range("H109").Value = 31500
# row 126
$ws.Range("H126").Value = 8894.111
$ws.Range("I126").Value = 3426.182
$ws.Range("K126").Value = 10278.546
$ws.Range("M126").Value = -7808.545999999998

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 34
$ws.Range("H34").Value = 1799.8572
$ws.Range("J34").Value = 2459.8
$ws.Range("L34").Value = 7379.400000000001
$ws.Range("N34").Value = -7547.400000000001
# row 47
$ws.Range("H47").Value = 1488.2858
$ws.Range("I47").Value = 139.33333
$ws.Range("K47").Value = 417.99999
$ws.Range("M47").Value = 13.00001000000003
# row 68
$ws.Range("H68").Value = 1851.7297
$ws.Range("I68").Value = 1199.8334
$ws.Range("J68").Value = 2296.2046
$ws.Range("K68").Value = 3599.5002
$ws.Range("L68").Value = 6888.6138
$ws.Range("M68").Value = -2788.5002
$ws.Range("N68").Value = -8510.6138
# row 71
$ws.Range("H71").Value = 1851.7297
$ws.Range("I71").Value = 1199.8334
$ws.Range("J71").Value = 2296.2046
$ws.Range("K71").Value = 10798.5006
$ws.Range("L71").Value = 20665.8414
$ws.Range("M71").Value = -6742.500599999999
$ws.Range("N71").Value = -28777.8414
# row 94
$ws.Range("H94").Value = 4854.1665
$ws.Range("I94").Value = 3400
$ws.Range("J94").Value = 5892.857
$ws.Range("K94").Value = 10200
$ws.Range("L94").Value = 17678.571
$ws.Range("M94").Value = -9524
$ws.Range("N94").Value = -19030.571
# row 107
$ws.Range("H107").Value = 263034.44
$ws.Range("I107").Value = 404.09375
$ws.Range("J107").Value = 663233.1
$ws.Range("K107").Value = 1212.28125
$ws.Range("L107").Value = 1989699.3
$ws.Range("M107").Value = 707.71875
$ws.Range("N107").Value = -1993539.3
# row 131
$ws.Range("H131").Value = 851.3418
$ws.Range("I131").Value = 703.3333
$ws.Range("J131").Value = 857.1842
$ws.Range("K131").Value = 2109.9999
$ws.Range("L131").Value = 2571.5526
$ws.Range("M131").Value = 2930.0001
$ws.Range("N131").Value = -12651.5526

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 59
$ws.Range("H59").Value = 5500
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 10000
$ws.Range("K59").Value = 1000
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = -417
$ws.Range("N59").Value = -11166

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 102036.8
$ws.Range("I40").Value = 144341.14
$ws.Range("J40").Value = 3326.6667
$ws.Range("K40").Value = 144341.14
$ws.Range("L40").Value = 3326.6667
$ws.Range("M40").Value = -144205.14
$ws.Range("N40").Value = -3598.6667
# row 100
$ws.Range("H100").Value = 2149.75
$ws.Range("I100").Value = 2200
$ws.Range("J100").Value = 2099.5
$ws.Range("K100").Value = 2200
$ws.Range("L100").Value = 2099.5
$ws.Range("M100").Value = -1659
$ws.Range("N100").Value = -3181.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 115
$ws.Range("H115").Value = 34991.332
$ws.Range("J115").Value = 34991.332
$ws.Range("L115").Value = 34991.332
$ws.Range("N115").Value = -38125.332
# row 122
$ws.Range("H122").Value = 1187.5555
$ws.Range("I122").Value = 1211
$ws.Range("K122").Value = 3633
$ws.Range("M122").Value = -1183
